$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-pulled dSF (column F) values differ from the stale dS0 copy for these rows.
$updates = @{
    7  = 5
    10 = -8
    12 = -2
    13 = 2
    19 = -1
    23 = -1
    27 = 3
    28 = -1
    33 = 0
    36 = -2
    37 = 1
    39 = 1
    42 = 1
    46 = -1
    47 = -2
    48 = -2
    52 = -2
    59 = -1
    62 = 4
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
